$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new content is purely digits (or a currency-style number) need
# to be forced to Text format first, otherwise the engine auto-coerces them
# to numeric/currency values and we lose the exact string the diff expects
# (e.g. "1023.78€" silently turning into the number 1023.78).
$ws.Range("C2").NumberFormat = "@"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("L2").NumberFormat = "@"
$ws.Range("K2:K5").NumberFormat = "@"

# --- Row 2: receiver / courier / product info updated ---
$ws.Range("A2").Value = "ΤΑΤΣΗΣ  ΓΕΩΡΓΙΟΣ"
$ws.Range("B2").Value = "ΕΛΤΑ COURIER"
$ws.Range("C2").Value = "134312928"
$ws.Range("D2").Value = "2681071591"
$ws.Range("E2").ClearContents()
$ws.Range("F2").Value = "47100"
$ws.Range("G2").Value = "ΑΡΤΑ"
$ws.Range("H2").Value = "ΦΛΕΜΙΝΓΚ ΚΑΙ ΠΕΡΙΦ ΟΔΟΣ 0"
$ws.Range("I2").Value = "GPT-0532"
$ws.Range("J2").Value = "OKI TONER M C3100/3000/3200/5100/5150/5200/5300/5400/5510 MAGENTA ΣΥΜΒΑΤΟ 3000 ΣΕΛΙΔΕΣ"
$ws.Range("K2").Value = "108"
$ws.Range("L2").Value = "1023.78€"

# --- Row 3: product info updated, receiver columns stay blank ---
$ws.Range("I3").Value = "GPI-0134"
$ws.Range("J3").Value = "HP INK No 88XL -  C9392A MAGENTA ΣΥΜΒΑΤΟ 28ml"
$ws.Range("K3").Value = "76"

# --- Row 4: product info updated, receiver columns stay blank ---
$ws.Range("I4").Value = "GPI-0023"
$ws.Range("J4").Value = "EPSON INK No 26XL - T2634XL YELLOW ΣΥΜΒΑΤΟ 10ml"
$ws.Range("K4").Value = "20"

# --- Row 5: brand-new row, only the product columns are populated ---
$ws.Range("I5").Value = "GPT-0142"
$ws.Range("J5").Value = "HP TONER CF383A MAGENTA ΣΥΜΒΑΤΟ 2700 ΣΕΛΙΔΕΣ"
$ws.Range("K5").Value = "78"
